# Fruta / hortaliza, semanal
# Insert a new weekly price record as row 190 (pushing existing rows 190-219
# down to 191-220) in the "Vega Modelo de Temuco - Perejil" dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 190; this shifts rows 190:219
# down to 191:220 and extends the used range to row 220.
$ws.Rows(190).Insert()

# Populate the newly inserted row 190 with the new weekly record.
$ws.Cells.Item(190, 1).Value = 10
$ws.Cells.Item(190, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(190, 3).Value = "La Araucanía"
$ws.Cells.Item(190, 4).Value = 44504
$ws.Cells.Item(190, 5).Value = 9
$ws.Cells.Item(190, 6).Value = 100112044
$ws.Cells.Item(190, 7).Value = "Perejil"
$ws.Cells.Item(190, 8).Value = "Sin especificar"
$ws.Cells.Item(190, 9).Value = "Primera"
$ws.Cells.Item(190, 10).Value = 45
$ws.Cells.Item(190, 11).Value = 5000
$ws.Cells.Item(190, 12).Value = 5000
$ws.Cells.Item(190, 13).Value = 5000
$ws.Cells.Item(190, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(190, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(190, 16).Value = 1667
$ws.Cells.Item(190, 17).Value = 3
$ws.Cells.Item(190, 18).Value = "Hortaliza"
